$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (PT Borwita): update lease end date, add 3-month reminder, decrease lease duration
$ws.Range("C5").Value = 45518
$ws.Range("D5").Value = "Reminder: Lease Ending Soon"
$ws.Range("E5").Value = 1

# Row 11 (PT Tumbakmas Niaga Sakti (Sasa)): update lease end date, add 3-month reminder, decrease lease duration
$ws.Range("C11").Value = 45657
$ws.Range("D11").Value = "Reminder: Lease Ending Soon"
$ws.Range("E11").Value = 1
